# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values computed for rows 2-9 (data rows), column G.
$newValues = @{
    2 = 0
    3 = 1
    4 = 4
    5 = 0
    6 = 2
    7 = 2
    8 = 0
    9 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
